$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6662.6665
$ws.Range("J32").Value = 6662.6665
$ws.Range("L32").Value = 6662.6665
$ws.Range("N32").Value = -7314.6665
$ws.Range("H43").Value = 4085.4167
$ws.Range("J43").Value = 4303.6
$ws.Range("L43").Value = 4303.6
$ws.Range("N43").Value = -4441.6
$ws.Range("H64").Value = 5991.25
$ws.Range("J64").Value = 5991.25
$ws.Range("L64").Value = 5991.25
$ws.Range("N64").Value = -6487.25
$ws.Range("H67").Value = 5991.25
$ws.Range("J67").Value = 5991.25
$ws.Range("L67").Value = 5991.25
$ws.Range("N67").Value = -7707.25
$ws.Range("H76").Value = 6688.048
$ws.Range("I76").Value = 6688.048
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6688.048
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6373.048
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 6688.048
$ws.Range("I79").Value = 6688.048
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6688.048
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5596.048
$ws.Range("N79").ClearContents()
$ws.Range("H107").Value = 1464.625
$ws.Range("I107").Value = 1464.625
$ws.Range("K107").Value = 1464.625
$ws.Range("M107").Value = 455.375
$ws.Range("H132").Value = 4325.8857
$ws.Range("I132").Value = 3571.0588
$ws.Range("K132").Value = 10713.1764
$ws.Range("M132").Value = -8183.1764
$ws.Range("H138").Value = 6050.107
$ws.Range("I138").Value = 2387
$ws.Range("J138").Value = 7271.143
$ws.Range("K138").Value = 7161
$ws.Range("L138").Value = 21813.429
$ws.Range("M138").Value = -2021
$ws.Range("N138").Value = -32093.429
$ws.Range("H141").Value = 6703.973
$ws.Range("I141").Value = 6092.6177
$ws.Range("K141").Value = 18277.8531
$ws.Range("M141").Value = -13097.8531

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 840.875
$ws.Range("I2").Value = 754.2143
$ws.Range("K2").Value = 754.2143
$ws.Range("M2").Value = -641.2143
$ws.Range("H63").Value = 2027.08
$ws.Range("I63").Value = 2041.125
$ws.Range("J63").Value = 1690
$ws.Range("K63").Value = 2041.125
$ws.Range("L63").Value = 1690
$ws.Range("M63").Value = -1355.125
$ws.Range("N63").Value = -3062
$ws.Range("H66").Value = 2027.08
$ws.Range("I66").Value = 2041.125
$ws.Range("J66").Value = 1690
$ws.Range("K66").Value = 10205.625
$ws.Range("L66").Value = 8450
$ws.Range("M66").Value = -6773.625
$ws.Range("N66").Value = -15314
$ws.Range("H110").Value = 3973
$ws.Range("I110").Value = 3463.3572
$ws.Range("K110").Value = 3463.3572
$ws.Range("M110").Value = -1418.3572
$ws.Range("H116").Value = 840.875
$ws.Range("I116").Value = 754.2143
$ws.Range("K116").Value = 754.2143
$ws.Range("M116").Value = 1539.7857
$ws.Range("H122").Value = 1893
$ws.Range("I122").Value = 1943
$ws.Range("J122").Value = 1543
$ws.Range("K122").Value = 5829
$ws.Range("L122").Value = 4629
$ws.Range("M122").Value = -3379
$ws.Range("N122").Value = -9529
$ws.Range("H132").Value = 21622.113
$ws.Range("I132").Value = 23876.127
$ws.Range("K132").Value = 71628.38099999999
$ws.Range("M132").Value = -69098.38099999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 840.875
$ws.Range("I3").Value = 754.2143
$ws.Range("K3").Value = 754.2143
$ws.Range("M3").Value = -640.2143
$ws.Range("H86").Value = 2615.2727
$ws.Range("I86").Value = 1682.1428
$ws.Range("J86").Value = 4248.25
$ws.Range("K86").Value = 1682.1428
$ws.Range("L86").Value = 4248.25
$ws.Range("M86").Value = -559.1428000000001
$ws.Range("N86").Value = -6494.25
$ws.Range("H89").Value = 2615.2727
$ws.Range("I89").Value = 1682.1428
$ws.Range("J89").Value = 4248.25
$ws.Range("K89").Value = 8410.714
$ws.Range("L89").Value = 21241.25
$ws.Range("M89").Value = -2794.714
$ws.Range("N89").Value = -32473.25
$ws.Range("H105").Value = 2998
$ws.Range("I105").Value = 2496.75
$ws.Range("K105").Value = 2496.75
$ws.Range("M105").Value = -749.75
$ws.Range("H107").Value = 2349.3572
$ws.Range("I107").Value = 2063.3635
$ws.Range("K107").Value = 2063.3635
$ws.Range("M107").Value = -143.3634999999999
$ws.Range("H134").Value = 2290.1086
$ws.Range("I134").Value = 2187.6135
$ws.Range("J134").Value = 4545
$ws.Range("K134").Value = 6562.8405
$ws.Range("L134").Value = 13635
$ws.Range("M134").Value = -4027.8405
$ws.Range("N134").Value = -18705

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.6875
$ws.Range("J7").Value = 109.875
$ws.Range("L7").Value = 109.875
$ws.Range("N7").Value = -335.875
$ws.Range("H22").Value = 550.8889
$ws.Range("I22").Value = 274
$ws.Range("K22").Value = 274
$ws.Range("M22").Value = 76
$ws.Range("H62").Value = 2857.1428
$ws.Range("I62").Value = 2401
$ws.Range("J62").Value = 3997.5
$ws.Range("K62").Value = 2401
$ws.Range("L62").Value = 3997.5
$ws.Range("M62").Value = -1777
$ws.Range("N62").Value = -5245.5
$ws.Range("H65").Value = 2857.1428
$ws.Range("I65").Value = 2401
$ws.Range("J65").Value = 3997.5
$ws.Range("K65").Value = 12005
$ws.Range("L65").Value = 19987.5
$ws.Range("M65").Value = -8885
$ws.Range("N65").Value = -26227.5
$ws.Range("H132").Value = 4257.5386
$ws.Range("I132").Value = 4347.1
$ws.Range("K132").Value = 13041.3
$ws.Range("M132").Value = -10511.3
$ws.Range("H134").Value = 31014.914
$ws.Range("I134").Value = 32705.213
$ws.Range("K134").Value = 98115.639
$ws.Range("M134").Value = -95580.639

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2074.7585
$ws.Range("J5").Value = 3008.3125
$ws.Range("L5").Value = 9024.9375
$ws.Range("N5").Value = -9248.9375
$ws.Range("H92").Value = 818.12
$ws.Range("J92").Value = 1136
$ws.Range("L92").Value = 3408
$ws.Range("N92").Value = -5904
$ws.Range("H135").Value = 2074.7585
$ws.Range("J135").Value = 3008.3125
$ws.Range("L135").Value = 27074.8125
$ws.Range("N135").Value = -32144.8125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5686.143
$ws.Range("I70").Value = 4898.5
$ws.Range("K70").Value = 4898.5
$ws.Range("M70").Value = -4628.5
$ws.Range("H73").Value = 5686.143
$ws.Range("I73").Value = 4898.5
$ws.Range("K73").Value = 4898.5
$ws.Range("M73").Value = -3962.5
$ws.Range("H80").Value = 3431.8
$ws.Range("J80").Value = 3817.8
$ws.Range("L80").Value = 3817.8
$ws.Range("N80").Value = -5813.8
$ws.Range("H83").Value = 3431.8
$ws.Range("J83").Value = 3817.8
$ws.Range("L83").Value = 19089
$ws.Range("N83").Value = -29073
$ws.Range("H132").Value = 32473.617
$ws.Range("I132").Value = 34584
$ws.Range("J132").Value = 10666.333
$ws.Range("K132").Value = 103752
$ws.Range("L132").Value = 31998.999
$ws.Range("M132").Value = -101222
$ws.Range("N132").Value = -37058.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5875.8887
$ws.Range("I68").Value = 4798
$ws.Range("K68").Value = 4798
$ws.Range("M68").Value = -4049
$ws.Range("H71").Value = 5875.8887
$ws.Range("I71").Value = 4798
$ws.Range("K71").Value = 23990
$ws.Range("M71").Value = -20246
$ws.Range("H82").Value = 3114.9546
$ws.Range("I82").Value = 1805.7778
$ws.Range("J82").Value = 4021.3076
$ws.Range("K82").Value = 1805.7778
$ws.Range("L82").Value = 4021.3076
$ws.Range("M82").Value = -1444.7778
$ws.Range("N82").Value = -4743.3076
$ws.Range("H85").Value = 3114.9546
$ws.Range("I85").Value = 1805.7778
$ws.Range("J85").Value = 4021.3076
$ws.Range("K85").Value = 1805.7778
$ws.Range("L85").Value = 4021.3076
$ws.Range("M85").Value = -557.7778000000001
$ws.Range("N85").Value = -6517.3076
$ws.Range("H132").Value = 32576.61
$ws.Range("I132").Value = 46791.035
$ws.Range("K132").Value = 140373.105
$ws.Range("M132").Value = -137843.105
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 631777.9399999999
$ws.Range("I81").Value = 1450.3636
$ws.Range("J81").Value = 2018498.6
$ws.Range("K81").Value = 2900.7272
$ws.Range("L81").Value = 4036997.2
$ws.Range("M81").Value = -1839.7272
$ws.Range("N81").Value = -4039119.2
$ws.Range("H84").Value = 631777.9399999999
$ws.Range("I84").Value = 1450.3636
$ws.Range("J84").Value = 2018498.6
$ws.Range("K84").Value = 14503.636
$ws.Range("L84").Value = 20184986
$ws.Range("M84").Value = -9199.635999999999
$ws.Range("N84").Value = -20195594
$ws.Range("H94").Value = 47775.8
$ws.Range("H100").Value = 1146.2142
$ws.Range("J100").Value = 1229.3334
$ws.Range("L100").Value = 2458.6668
$ws.Range("N100").Value = -3540.6668
$ws.Range("H122").Value = 3540
$ws.Range("I122").Value = 3122.1428
$ws.Range("K122").Value = 9366.428400000001
$ws.Range("M122").Value = -6916.428400000001
$ws.Range("H136").Value = 5167.478
$ws.Range("I136").Value = 5319.067
$ws.Range("J136").Value = 4883.25
$ws.Range("K136").Value = 15957.201
$ws.Range("L136").Value = 14649.75
$ws.Range("M136").Value = -13407.201
$ws.Range("N136").Value = -19749.75
